$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Force Text number format on cells whose new numeric-looking value would
# --- otherwise be auto-converted to a Number by Excel (column D price values).
$textForceRange = $ws.Range("D4","D5","D6","D10","D11","D12","D13","D15","D16","D21","D22","D24","D25","D26","D27","D28","D30","D31","D32","D34","D35","D36","D39","D41","D42","D43","D45","D47","D48","D51","D37","D38","D49","D50")
$textForceRange.NumberFormat = "@"

# --- Column D / E value updates ---
$ws.Range("D2").Value = "43.260.66"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "2.282.47"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "114.84"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "304.42"
$ws.Range("E6").Value = "  +7.38%  "
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "45.04"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "55.11"
$ws.Range("D13").Value = "8.91"
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("E14").Value = "  +19.88%  "
$ws.Range("D15").Value = "0.105"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "15.44"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "2.620.35"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").Value = "2.283.36"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "43.189.86"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "7.25"
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("D22").Value = "74.88"
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("E23").Value = "  +12.02%  "
$ws.Range("D24").Value = "2.47"
$ws.Range("E24").Value = "  +4.97%  "
$ws.Range("D25").Value = "254.65"
$ws.Range("E25").Value = "  +10.12%  "
$ws.Range("D26").Value = "9.05"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "11.74"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").Value = "38.35"
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("D31").Value = "175.43"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "22.12"
$ws.Range("E32").Value = "  +4.64%  "
$ws.Range("E33").Value = "  -3.58%  "
$ws.Range("D34").Value = "0.0901"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "5.73"
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("D36").Value = "5.07"
$ws.Range("E36").Value = "  +9.26%  "
$ws.Range("D39").Value = "0.0378"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").Value = "2.54"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("D42").Value = "72.87"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "0.233"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "12.67"
$ws.Range("E45").Value = "  -5.63%  "
$ws.Range("E46").Value = "  +3.02%  "
$ws.Range("D47").Value = "5.64"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("D48").Value = "107.43"
$ws.Range("E48").Value = "  +6.22%  "
$ws.Range("D51").Value = "73.88"
$ws.Range("E51").Value = "  +6.28%  "

# --- Row swaps (37/38 and 49/50): coin + link + price + volume all change ---
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "4.30"
$ws.Range("E37").Value = "  -7.36%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "0.129"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "8.80"
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "1.30"
$ws.Range("E50").Value = "  +0.52%  "
